$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "Largus Фургон New" entry (currently row 127, id 454) so the
#     list stays sorted alphabetically by model name after it is renamed
#     from "Фургон New" to "Largus Фургон New".
#     Copy row 127 and insert it before row 110, then remove the old row
#     (which, after the insert, has shifted down to row 128).
$ws.Rows(127).Copy()
$ws.Rows(110).Insert()
$ws.Rows(128).Delete()

# Rename the model text of the newly inserted row 110.
$ws.Range("C110").Value = "Largus Фургон New"

# --- Move the "Vesta CNG" entry (now at row 116, id 436) so the list stays
#     sorted alphabetically after it is renamed to "Vesta Sedan CNG" (which
#     sorts right before "Vesta Sedan New").
#     Copy row 116 and insert it before row 124, then remove the original
#     row 116 (its position is unaffected by an insert further down).
$ws.Rows(116).Copy()
$ws.Rows(124).Insert()
$ws.Rows(116).Delete()

# Rename the model text of the relocated row (now at row 123).
$ws.Range("C123").Value = "Vesta Sedan CNG"

Write-Host "done"
